$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 'M1'
$ws.Cells.Item(2,2).Value = 'Il1a'
$ws.Cells.Item(2,3).Value = 'Il1r1'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.056073
$ws.Cells.Item(2,8).Value = 9.168219
$ws.Cells.Item(2,9).Value = 0.4559514113020136
$ws.Cells.Item(2,10).Value = 0.4559514113020135
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 20.003843
$ws.Cells.Item(2,14).Value = 40.007686
$ws.Cells.Item(2,15).Value = 0.1517504313331627
$ws.Cells.Item(2,16).Value = 0.1111693208434551
$ws.Cells.Item(2,17).Value = 61.133204488539
$ws.Cells.Item(2,18).Value = 366.799226931234
$ws.Cells.Item(2,19).Value = 0.06919082333204485
$ws.Cells.Item(2,20).Value = 0.05068780873205971

$ws.Cells.Item(3,1).Value = 'M1'
$ws.Cells.Item(3,2).Value = 'Il1a'
$ws.Cells.Item(3,3).Value = 'Il1r1'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 3.056073
$ws.Cells.Item(3,8).Value = 9.168219
$ws.Cells.Item(3,9).Value = 0.4559514113020136
$ws.Cells.Item(3,10).Value = 0.4559514113020135
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 85.28390633333333
$ws.Cells.Item(3,14).Value = 255.851719
$ws.Cells.Item(3,15).Value = 0.6469691634682588
$ws.Cells.Item(3,16).Value = 0.710934939797831
$ws.Cells.Item(3,17).Value = 260.633843479829
$ws.Cells.Item(3,18).Value = 2345.704591318461
$ws.Cells.Item(3,19).Value = 0.2949865031522357
$ws.Cells.Item(3,20).Value = 0.324151789144733

$ws.Cells.Item(4,1).Value = 'M1'
$ws.Cells.Item(4,2).Value = 'Il1a'
$ws.Cells.Item(4,3).Value = 'Il1r1'
$ws.Cells.Item(4,4).Value = 'M1'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 3.056073
$ws.Cells.Item(4,8).Value = 9.168219
$ws.Cells.Item(4,9).Value = 0.4559514113020136
$ws.Cells.Item(4,10).Value = 0.4559514113020135
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.3719913333333333
$ws.Cells.Item(4,14).Value = 1.115974
$ws.Cells.Item(4,15).Value = 0.002821950026578976
$ws.Cells.Item(4,16).Value = 0.003100955942789443
$ws.Cells.Item(4,17).Value = 1.136832670034
$ws.Cells.Item(4,18).Value = 10.231494030306
$ws.Cells.Item(4,19).Value = 0.001286672097242439
$ws.Cells.Item(4,20).Value = 0.001413885238500212

$ws.Cells.Item(5,1).Value = 'M1'
$ws.Cells.Item(5,2).Value = 'Il1a'
$ws.Cells.Item(5,3).Value = 'Il1r1'
$ws.Cells.Item(5,4).Value = 'M2'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3.056073
$ws.Cells.Item(5,8).Value = 9.168219
$ws.Cells.Item(5,9).Value = 0.4559514113020136
$ws.Cells.Item(5,10).Value = 0.4559514113020135
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.3862156666666667
$ws.Cells.Item(5,14).Value = 1.158647
$ws.Cells.Item(5,15).Value = 0.002929856728244252
$ws.Cells.Item(5,16).Value = 0.003219531369230071
$ws.Cells.Item(5,17).Value = 1.180303271077
$ws.Cells.Item(5,18).Value = 10.622729439693
$ws.Cells.Item(5,19).Value = 0.001335872310155666
$ws.Cells.Item(5,20).Value = 0.001467949871531555

$ws.Cells.Item(6,1).Value = 'M1'
$ws.Cells.Item(6,2).Value = 'Il1a'
$ws.Cells.Item(6,3).Value = 'Il1r1'
$ws.Cells.Item(6,4).Value = 'Neutro'
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3.056073
$ws.Cells.Item(6,8).Value = 9.168219
$ws.Cells.Item(6,9).Value = 0.4559514113020136
$ws.Cells.Item(6,10).Value = 0.4559514113020135
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 10.19719066666667
$ws.Cells.Item(6,14).Value = 30.591572
$ws.Cells.Item(6,15).Value = 0.07735654004348905
$ws.Cells.Item(6,16).Value = 0.08500477340213224
$ws.Cells.Item(6,17).Value = 31.163359072252
$ws.Cells.Item(6,18).Value = 280.470231650268
$ws.Cells.Item(6,19).Value = 0.03527082360626956
$ws.Cells.Item(6,20).Value = 0.03875804640011005

$ws.Cells.Item(7,1).Value = 'M1'
$ws.Cells.Item(7,2).Value = 'Il1a'
$ws.Cells.Item(7,3).Value = 'Il1r1'
$ws.Cells.Item(7,4).Value = 'sCs'
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 3.056073
$ws.Cells.Item(7,8).Value = 9.168219
$ws.Cells.Item(7,9).Value = 0.4559514113020136
$ws.Cells.Item(7,10).Value = 0.4559514113020135
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 15.5775195
$ws.Cells.Item(7,14).Value = 31.155039
$ws.Cells.Item(7,15).Value = 0.1181720584002661
$ws.Cells.Item(7,16).Value = 0.08657047864456238
$ws.Cells.Item(7,17).Value = 47.6060367509235
$ws.Cells.Item(7,18).Value = 285.636220505541
$ws.Cells.Item(7,19).Value = 0.05388071680406529
$ws.Cells.Item(7,20).Value = 0.03947193191507904

$ws.Cells.Item(8,1).Value = 'M2'
$ws.Cells.Item(8,2).Value = 'Il1a'
$ws.Cells.Item(8,3).Value = 'Il1r1'
$ws.Cells.Item(8,4).Value = 'ECs'
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 3.592139
$ws.Cells.Item(8,8).Value = 10.776417
$ws.Cells.Item(8,9).Value = 0.5359298834298145
$ws.Cells.Item(8,10).Value = 0.5359298834298145
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 20.003843
$ws.Cells.Item(8,14).Value = 40.007686
$ws.Cells.Item(8,15).Value = 0.1517504313331627
$ws.Cells.Item(8,16).Value = 0.1111693208434551
$ws.Cells.Item(8,17).Value = 71.856584590177
$ws.Cells.Item(8,18).Value = 431.139507541062
$ws.Cells.Item(8,19).Value = 0.08132759097480598
$ws.Cells.Item(8,20).Value = 0.05957896116060456

$ws.Cells.Item(9,1).Value = 'M2'
$ws.Cells.Item(9,2).Value = 'Il1a'
$ws.Cells.Item(9,3).Value = 'Il1r1'
$ws.Cells.Item(9,4).Value = 'FAPs'
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 3.592139
$ws.Cells.Item(9,8).Value = 10.776417
$ws.Cells.Item(9,9).Value = 0.5359298834298145
$ws.Cells.Item(9,10).Value = 0.5359298834298145
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 85.28390633333333
$ws.Cells.Item(9,14).Value = 255.851719
$ws.Cells.Item(9,15).Value = 0.6469691634682588
$ws.Cells.Item(9,16).Value = 0.710934939797831
$ws.Cells.Item(9,17).Value = 306.3516460123137
$ws.Cells.Item(9,18).Value = 2757.164814110823
$ws.Cells.Item(9,19).Value = 0.3467301083602286
$ws.Cells.Item(9,20).Value = 0.3810112794120337

$ws.Cells.Item(10,1).Value = 'M2'
$ws.Cells.Item(10,2).Value = 'Il1a'
$ws.Cells.Item(10,3).Value = 'Il1r1'
$ws.Cells.Item(10,4).Value = 'M1'
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 3.592139
$ws.Cells.Item(10,8).Value = 10.776417
$ws.Cells.Item(10,9).Value = 0.5359298834298145
$ws.Cells.Item(10,10).Value = 0.5359298834298145
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.3719913333333333
$ws.Cells.Item(10,14).Value = 1.115974
$ws.Cells.Item(10,15).Value = 0.002821950026578976
$ws.Cells.Item(10,16).Value = 0.003100955942789443
$ws.Cells.Item(10,17).Value = 1.336244576128667
$ws.Cells.Item(10,18).Value = 12.026201185158
$ws.Cells.Item(10,19).Value = 0.001512367348789232
$ws.Cells.Item(10,20).Value = 0.001661894956940137

$ws.Cells.Item(11,1).Value = 'M2'
$ws.Cells.Item(11,2).Value = 'Il1a'
$ws.Cells.Item(11,3).Value = 'Il1r1'
$ws.Cells.Item(11,4).Value = 'M2'
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 3.592139
$ws.Cells.Item(11,8).Value = 10.776417
$ws.Cells.Item(11,9).Value = 0.5359298834298145
$ws.Cells.Item(11,10).Value = 0.5359298834298145
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.3862156666666667
$ws.Cells.Item(11,14).Value = 1.158647
$ws.Cells.Item(11,15).Value = 0.002929856728244252
$ws.Cells.Item(11,16).Value = 0.003219531369230071
$ws.Cells.Item(11,17).Value = 1.387340358644333
$ws.Cells.Item(11,18).Value = 12.486063227799
$ws.Cells.Item(11,19).Value = 0.001570197774833999
$ws.Cells.Item(11,20).Value = 0.001725443071410103

$ws.Cells.Item(12,1).Value = 'M2'
$ws.Cells.Item(12,2).Value = 'Il1a'
$ws.Cells.Item(12,3).Value = 'Il1r1'
$ws.Cells.Item(12,4).Value = 'Neutro'
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 3.592139
$ws.Cells.Item(12,8).Value = 10.776417
$ws.Cells.Item(12,9).Value = 0.5359298834298145
$ws.Cells.Item(12,10).Value = 0.5359298834298145
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 10.19719066666667
$ws.Cells.Item(12,14).Value = 30.591572
$ws.Cells.Item(12,15).Value = 0.07735654004348905
$ws.Cells.Item(12,16).Value = 0.08500477340213224
$ws.Cells.Item(12,17).Value = 36.62972628416933
$ws.Cells.Item(12,18).Value = 329.667536557524
$ws.Cells.Item(12,19).Value = 0.04145768148804086
$ws.Cells.Item(12,20).Value = 0.04555659830038253

$ws.Cells.Item(13,1).Value = 'M2'
$ws.Cells.Item(13,2).Value = 'Il1a'
$ws.Cells.Item(13,3).Value = 'Il1r1'
$ws.Cells.Item(13,4).Value = 'sCs'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 3.592139
$ws.Cells.Item(13,8).Value = 10.776417
$ws.Cells.Item(13,9).Value = 0.5359298834298145
$ws.Cells.Item(13,10).Value = 0.5359298834298145
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 15.5775195
$ws.Cells.Item(13,14).Value = 31.155039
$ws.Cells.Item(13,15).Value = 0.1181720584002661
$ws.Cells.Item(13,16).Value = 0.08657047864456238
$ws.Cells.Item(13,17).Value = 55.95661531921051
$ws.Cells.Item(13,18).Value = 335.7396919152631
$ws.Cells.Item(13,19).Value = 0.06333193748311583
$ws.Cells.Item(13,20).Value = 0.04639570652844356

$ws.Cells.Item(14,1).Value = 'Neutro'
$ws.Cells.Item(14,2).Value = 'Il1a'
$ws.Cells.Item(14,3).Value = 'Il1r1'
$ws.Cells.Item(14,4).Value = 'ECs'
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.05441666666666667
$ws.Cells.Item(14,8).Value = 0.16325
$ws.Cells.Item(14,9).Value = 0.0081187052681719
$ws.Cells.Item(14,10).Value = 0.008118705268171898
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 20.003843
$ws.Cells.Item(14,14).Value = 40.007686
$ws.Cells.Item(14,15).Value = 0.1517504313331627
$ws.Cells.Item(14,16).Value = 0.1111693208434551
$ws.Cells.Item(14,17).Value = 1.088542456583333
$ws.Cells.Item(14,18).Value = 6.5312547395
$ws.Cells.Item(14,19).Value = 0.001232017026311906
$ws.Cells.Item(14,20).Value = 0.0009025509507908511

$ws.Cells.Item(15,1).Value = 'Neutro'
$ws.Cells.Item(15,2).Value = 'Il1a'
$ws.Cells.Item(15,3).Value = 'Il1r1'
$ws.Cells.Item(15,4).Value = 'FAPs'
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.05441666666666667
$ws.Cells.Item(15,8).Value = 0.16325
$ws.Cells.Item(15,9).Value = 0.0081187052681719
$ws.Cells.Item(15,10).Value = 0.008118705268171898
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 85.28390633333333
$ws.Cells.Item(15,14).Value = 255.851719
$ws.Cells.Item(15,15).Value = 0.6469691634682588
$ws.Cells.Item(15,16).Value = 0.710934939797831
$ws.Cells.Item(15,17).Value = 4.640865902972222
$ws.Cells.Item(15,18).Value = 41.76779312675
$ws.Cells.Item(15,19).Value = 0.00525255195579452
$ws.Cells.Item(15,20).Value = 0.005771871241064121

$ws.Cells.Item(16,1).Value = 'Neutro'
$ws.Cells.Item(16,2).Value = 'Il1a'
$ws.Cells.Item(16,3).Value = 'Il1r1'
$ws.Cells.Item(16,4).Value = 'M1'
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.05441666666666667
$ws.Cells.Item(16,8).Value = 0.16325
$ws.Cells.Item(16,9).Value = 0.0081187052681719
$ws.Cells.Item(16,10).Value = 0.008118705268171898
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.3719913333333333
$ws.Cells.Item(16,14).Value = 1.115974
$ws.Cells.Item(16,15).Value = 0.002821950026578976
$ws.Cells.Item(16,16).Value = 0.003100955942789443
$ws.Cells.Item(16,17).Value = 0.02024252838888889
$ws.Cells.Item(16,18).Value = 0.1821827555
$ws.Cells.Item(16,19).Value = 0.00002291058054730456
$ws.Cells.Item(16,20).Value = 0.0000251757473490936

$ws.Cells.Item(17,1).Value = 'Neutro'
$ws.Cells.Item(17,2).Value = 'Il1a'
$ws.Cells.Item(17,3).Value = 'Il1r1'
$ws.Cells.Item(17,4).Value = 'M2'
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.05441666666666667
$ws.Cells.Item(17,8).Value = 0.16325
$ws.Cells.Item(17,9).Value = 0.0081187052681719
$ws.Cells.Item(17,10).Value = 0.008118705268171898
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.3862156666666667
$ws.Cells.Item(17,14).Value = 1.158647
$ws.Cells.Item(17,15).Value = 0.002929856728244252
$ws.Cells.Item(17,16).Value = 0.003219531369230071
$ws.Cells.Item(17,17).Value = 0.02101656919444445
$ws.Cells.Item(17,18).Value = 0.18914912275
$ws.Cells.Item(17,19).Value = 0.00002378664325458549
$ws.Cells.Item(17,20).Value = 0.00002613842628841286

$ws.Cells.Item(18,1).Value = 'Neutro'
$ws.Cells.Item(18,2).Value = 'Il1a'
$ws.Cells.Item(18,3).Value = 'Il1r1'
$ws.Cells.Item(18,4).Value = 'Neutro'
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = 0.3333333333333333
$ws.Cells.Item(18,7).Value = 0.05441666666666667
$ws.Cells.Item(18,8).Value = 0.16325
$ws.Cells.Item(18,9).Value = 0.0081187052681719
$ws.Cells.Item(18,10).Value = 0.008118705268171898
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 10.19719066666667
$ws.Cells.Item(18,14).Value = 30.591572
$ws.Cells.Item(18,15).Value = 0.07735654004348905
$ws.Cells.Item(18,16).Value = 0.08500477340213224
$ws.Cells.Item(18,17).Value = 0.5548971254444445
$ws.Cells.Item(18,18).Value = 4.994074129
$ws.Cells.Item(18,19).Value = 0.0006280349491786251
$ws.Cells.Item(18,20).Value = 0.0006901287016396494

$ws.Cells.Item(19,1).Value = 'Neutro'
$ws.Cells.Item(19,2).Value = 'Il1a'
$ws.Cells.Item(19,3).Value = 'Il1r1'
$ws.Cells.Item(19,4).Value = 'sCs'
$ws.Cells.Item(19,5).Value = 1
$ws.Cells.Item(19,6).Value = 0.3333333333333333
$ws.Cells.Item(19,7).Value = 0.05441666666666667
$ws.Cells.Item(19,8).Value = 0.16325
$ws.Cells.Item(19,9).Value = 0.0081187052681719
$ws.Cells.Item(19,10).Value = 0.008118705268171898
$ws.Cells.Item(19,11).Value = 2
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 15.5775195
$ws.Cells.Item(19,14).Value = 31.155039
$ws.Cells.Item(19,15).Value = 0.1181720584002661
$ws.Cells.Item(19,16).Value = 0.08657047864456238
$ws.Cells.Item(19,17).Value = 0.8476766861250001
$ws.Cells.Item(19,18).Value = 5.086060116750001
$ws.Cells.Item(19,19).Value = 0.0009594041130849577
$ws.Cells.Item(19,20).Value = 0.0007028402010397715
